$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 23:38"

# --- Update Estados Unidos (row 4) totals ---
$ws.Range("B4").Value = 1156257
$ws.Range("C4").Value = 25227
$ws.Range("D4").Value = 170201
$ws.Range("E4").Value = 918831
$ws.Range("F4").Value = 16455
$ws.Range("G4").Value = 1472
$ws.Range("H4").Value = 67225

# --- Update Alemania (row 9) totals ---
$ws.Range("B9").Value = 164967
$ws.Range("C9").Value = 890
$ws.Range("D9").Value = 129000
$ws.Range("E9").Value = 29173
$ws.Range("F9").Value = 2105
$ws.Range("G9").Value = 58
$ws.Range("H9").Value = 6794

# --- Gabon overtakes Paraguay/Mauricio/Montenegro/Isla de Man/Guinea Ecuatorial ---
# Insert Gabon's updated row right after Venezuela (row 124), pushing the
# other five countries down one rank (their data is unchanged, just shifted).

$ws.Range("A125").Value = "Gabon"
$ws.Range("B125").Value = 335
$ws.Range("C125").Value = 59
$ws.Range("D125").Value = 85
$ws.Range("E125").Value = 245
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 5

$ws.Range("A126").Value = "Paraguay"
$ws.Range("B126").Value = 333
$ws.Range("C126").Value = 67
$ws.Range("D126").Value = 115
$ws.Range("E126").Value = 208
$ws.Range("F126").Value = 5
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 10

$ws.Range("A127").Value = "Mauricio"
$ws.Range("B127").Value = 332
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 314
$ws.Range("E127").Value = 8
$ws.Range("F127").Value = 3
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 10

$ws.Range("A128").Value = "Montenegro"
$ws.Range("B128").Value = 322
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 245
$ws.Range("E128").Value = 69
$ws.Range("F128").Value = 2
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 8

$ws.Range("A129").Value = "Isla de Man"
$ws.Range("B129").Value = 320
$ws.Range("C129").Value = 4
$ws.Range("D129").Value = 271
$ws.Range("E129").Value = 27
$ws.Range("F129").Value = 21
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 22

$ws.Range("A130").Value = "Guinea Ecuatorial"
$ws.Range("B130").Value = 315
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 9
$ws.Range("E130").Value = 305
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 1

# Row 131 (Vietnam) keeps its original data, nothing to change.
